$wb = $excel.ActiveWorkbook

# --- Quote_csv: delete the "Subject" text column (old column D) ---
# Columns E/F/G (Subject_1_ID/Subject_2_ID/Subject_3_ID) shift left to D/E/F.
$wsQuoteCsv = $wb.Worksheets.Item("Quote_csv")
$wsQuoteCsv.Columns("D").Delete()

# --- Update view/selection state on several sheets to match the saved session ---

# Subject sheet: drop the frozen top-left scroll position, move selection to C13
$wsSubject = $wb.Worksheets.Item("Subject")
$wsSubject.Activate()
$wsSubject.Range("C13").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Country sheet: move selection to E21
$wsCountry = $wb.Worksheets.Item("Country")
$wsCountry.Activate()
$wsCountry.Range("E21").Select()

# Quote_csv sheet: move selection to L32
$wsQuoteCsv.Activate()
$wsQuoteCsv.Range("L32").Select()

# Author_csv sheet: select the whole sheet (Ctrl+A) instead of a single cell
$wsAuthorCsv = $wb.Worksheets.Item("Author_csv")
$wsAuthorCsv.Activate()
$wsAuthorCsv.Cells.Select()

# Jobs sheet: move selection to D19 and leave this as the final active sheet/tab
$wsJobs = $wb.Worksheets.Item("Jobs")
$wsJobs.Activate()
$wsJobs.Range("D19").Select()
